# Update NATMI LR-pair data with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",   "Nxph3", "Nrxn1", "FAPs",  1, 0.3333333333333333, 0.07940033333333334, 0.238201,   0.2424709127739492, 0.3243800105402634, 2, 0.6666666666666666, 0.03298066666666667, 0.098942, 0.680931013599075, 0.7619715055833655, 0.002618675926888889, 0.023568083342,    0.1651059644034582, 0.2471683250125125),
    @("ECs",   "Nxph3", "Nrxn1", "MuSCs", 1, 0.3333333333333333, 0.07940033333333334, 0.238201,   0.2424709127739492, 0.3243800105402634, 1, 0.5,                0.015454,             0.030908, 0.319068986400925, 0.2380284944166346, 0.001227052751333333, 0.007362316508,    0.07736494837049107, 0.07721168552775096),
    @("MuSCs", "Nxph3", "Nrxn1", "FAPs",  1, 0.5,                0.248063,             0.496126,   0.7575290872260508, 0.6756199894597367, 2, 0.6666666666666666, 0.03298066666666667, 0.098942, 0.680931013599075, 0.7619715055833655, 0.008181283115333332, 0.049087698692,    0.5158250491956169, 0.5148031805708531),
    @("MuSCs", "Nxph3", "Nrxn1", "MuSCs", 1, 0.5,                0.248063,             0.496126,   0.7575290872260508, 0.6756199894597367, 1, 0.5,                0.015454,             0.030908, 0.319068986400925, 0.2380284944166346, 0.003833565602,       0.015334262408,    0.2417040380304339, 0.1608168088888837)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}
